# add code for cham cong function
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (column D) for rows 4, 5, 6, 8, 12, 13
$rows = @(4, 5, 6, 8, 12, 13)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "2024-07-06T13:10:00.000Z"
}

# Update properties.Số lượng đơn.number (AK8) from 58 to 59
$ws.Range("AK8").Value = 59
